$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range('D2').Value = '28.830.49'
$ws.Range('E2').Value = '  -1.29%  '

$ws.Range('D3').Value = '1.814.07'
$ws.Range('E3').Value = '  -0.81%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.69'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5918'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.29%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2758'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.21%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06756'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -4.74%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.91'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.16%  '

$ws.Range('E11').Value = '  -1.82%  '

$ws.Range('D12').Value = '1.810.67'
$ws.Range('E12').Value = '  -1.78%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.683'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.76%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6257'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.20%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009292'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -6.51%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '74.75'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -5.73%  '

$ws.Range('D17').Value = '28.637.60'
$ws.Range('E17').Value = '  -1.94%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.469'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -8.29%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.25%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '209.09'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -8.56%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.41'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.48%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.782'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.78%  '

$ws.Range('E23').Value = '  +0.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '154.57'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.58%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.818'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.65%  '

$ws.Range('E26').Value = '  -2.46%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.34'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.35%  '

$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06324'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.49%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.409'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.81%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.429'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.08%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.743'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.712'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.33%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.699'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.22%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.054'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.92%  '

$ws.Range('E35').Value = '  -3.29%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.522'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.31%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.733'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.95%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.447'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.07%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01695'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.95%  '

$ws.Range('D40').Value = '1.132.93'
$ws.Range('E40').Value = '  -8.21%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8670'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -6.64%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.004'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.27%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.44'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.56%  '

$ws.Range('E44').Value = '  -0.78%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.50'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.89%  '

$ws.Range('E46').Value = '  -3.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.576'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.25%  '


$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05451'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.87%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.272'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.54%  '

$ws.Range('E51').Value = '  +0.24%  '
